$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "bazwilks changes" paragraph: originally split across two runs with a
#    <w:proofErr spellStart/spellEnd> pair bracketing "bazwilks" (leftover
#    spell-check markup). Target is a single run "bazwilks changes" with no
#    proofErr markup at all.
#
#    We insert a brand-new (plain, unformatted) paragraph right after the
#    tainted one, give it the merged text, then delete the original
#    paragraph (proofErr element included) outright.
# ---------------------------------------------------------------------------

$bazParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($paraText -eq "bazwilks changes") {
        $bazParaIndex = $i
        break
    }
}

if ($bazParaIndex -gt 0) {
    $nextPara = $d.Paragraphs($bazParaIndex + 1)
    $nextPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs($bazParaIndex + 1)
    $newRange = $newPara.Range
    [void]$newRange.MoveEnd(1, -1)
    $newRange.InsertAfter("bazwilks changes")

    $oldPara = $d.Paragraphs($bazParaIndex)
    $oldPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) After the "Test – 15.03.2021" paragraph there used to be a single empty
#    trailing paragraph. It becomes three paragraphs of body text:
#    "Another test", "Test", "Test".
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = "Another test"

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = "Test"

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count).Range.Text = "Test"
